$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B/C/E columns (coin names, links, volume percentages) - plain text values
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').Value = '  +6.03%  '
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -2.10%  '
$ws.Range('E10').Value = '  +7.94%  '
$ws.Range('E11').Value = '  +11.43%  '
$ws.Range('E12').Value = '  -1.91%  '
$ws.Range('E13').Value = '  +2.07%  '
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').Value = '  -2.43%  '
$ws.Range('E20').Value = '  +1.13%  '
$ws.Range('E21').Value = '  +2.13%  '
$ws.Range('E22').Value = '  +5.20%  '
$ws.Range('E23').Value = '  -2.47%  '
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('B25').Value = 'Filecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E25').Value = '  +11.22%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('E28').Value = '  +2.74%  '
$ws.Range('E29').Value = '  +7.37%  '
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('E31').Value = '  -3.64%  '
$ws.Range('E32').Value = '  -3.74%  '
$ws.Range('E33').Value = '  +3.44%  '
$ws.Range('E34').Value = '  +33.12%  '
$ws.Range('E35').Value = '  -2.95%  '
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('E37').Value = '  -6.75%  '
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('E40').Value = '  +8.62%  '
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('E42').Value = '  +10.94%  '
$ws.Range('E43').Value = '  +0.69%  '
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('E47').Value = '  -0.48%  '
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('E51').Value = '  -1.23%  '

# Update D column (price values) - force text storage via Text format, then restore default style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.645.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.967.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '482.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.731'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.168'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000358'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.98'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.61'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.571.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.82'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.990.34'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.647.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '438.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.66'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.57'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.62'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '717.90'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.35'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.84'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0918'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '42.20'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '59.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.56'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0472'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.02'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.348'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.25'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '146.43'
$ws.Range('D50').Style = 'Normal'
